$d = $word.ActiveDocument
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$insertPoint = $d.Range($r.End, $r.End)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr><w:lastRenderedPageBreak/><w:t>Milestone 7</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">95 % </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Confidence interval for </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Distance between Airports (in miles)</w:t></w:r><w:r><w:t xml:space="preserve"> = (1187 – 1194)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>95% Confidence Interval for Average Fare Prices (in USD)</w:t></w:r><w:r><w:t xml:space="preserve"> = (218.70 – 219.31)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Since we used a 95% confidence level </w:t></w:r><w:r><w:t xml:space="preserve">the intervals that were computed come very close to each </w:t></w:r><w:r><w:t>other,</w:t></w:r><w:r><w:t xml:space="preserve"> especially with a large dataset. For the first interval it calculated the distance between airports to be from 1187 to 1194 miles.</w:t></w:r><w:r><w:t xml:space="preserve"> This means that we are 95% confident that true mean distance between these airlines lie somewhere between that interval</w:t></w:r><w:r><w:t xml:space="preserve"> and the same conclusion works for the second interval of Average Fare Price with intervals of $218.70 – $219.31. Since the sample size n is very large and our intervals aren’t extreme values away from each other we can say that these intervals can be trusted with more confidence.</w:t></w:r></w:p>'
$insertPoint.InsertXML($xml)
